# Applies the coinranking price/volume(1h) refresh from the Apr 15 2024
# GitHub Actions run: updated Price (D) / Volume(1h) (E) figures for every
# coin row, plus the Kaspa<->Maker (rows 40/41) and PEPE<->CoreDAO
# (rows 50/51) rank swaps that came with that run's re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, even when the text
# parses as a clean number (e.g. '554.20', '1.00', '0.120') — mirrors the
# source workbook, where every Price/Coin/Link cell is stored as a string.
# NumberFormat '@' forces the assignment to stick as text instead of being
# normalized to a number; ClearFormats() afterwards drops the temporary
# formatting again so the cell's style is left exactly as it started.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "63.068.68"
$ws.Range("E2").Value = "  -3.14%  "

$ws.Range("D3").Value = "3.100.50"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("E4").Value = "  +0.34%  "

Set-TextValue "D5" "554.20"
$ws.Range("E5").Value = "  -1.14%  "

Set-TextValue "D6" "137.48"
$ws.Range("E6").Value = "  -6.86%  "

Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("D8").Value = "3.093.82"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("E9").Value = "  -1.23%  "

Set-TextValue "D10" "6.64"
$ws.Range("E10").Value = "  -1.13%  "

Set-TextValue "D11" "0.159"
$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("E12").Value = "  -1.43%  "

Set-TextValue "D13" "35.04"
$ws.Range("E13").Value = "  -4.58%  "

$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "3.602.89"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "63.161.66"
$ws.Range("E16").Value = "  -2.66%  "

$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "3.104.05"
$ws.Range("E18").Value = "  +0.31%  "

Set-TextValue "D19" "503.71"
$ws.Range("E19").Value = "  -0.56%  "

Set-TextValue "D20" "6.67"
$ws.Range("E20").Value = "  -1.38%  "

Set-TextValue "D21" "13.53"
$ws.Range("E21").Value = "  -2.41%  "

Set-TextValue "D22" "0.706"
$ws.Range("E22").Value = "  +1.22%  "

Set-TextValue "D23" "7.25"
$ws.Range("E23").Value = "  -1.90%  "

Set-TextValue "D24" "77.60"
$ws.Range("E24").Value = "  -2.04%  "

Set-TextValue "D25" "12.29"
$ws.Range("E25").Value = "  -3.07%  "

Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.33%  "

Set-TextValue "D27" "2.75"
$ws.Range("E27").Value = "  -0.42%  "

Set-TextValue "D28" "8.27"
$ws.Range("E28").Value = "  -2.30%  "

Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.67%  "

Set-TextValue "D30" "1.95"
$ws.Range("E30").Value = "  -7.59%  "

Set-TextValue "D31" "26.25"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  -6.45%  "

$ws.Range("E33").Value = "  -2.73%  "

Set-TextValue "D34" "58.96"
$ws.Range("E34").Value = "  +11.81%  "

Set-TextValue "D35" "528.01"
$ws.Range("E35").Value = "  -10.51%  "

Set-TextValue "D36" "5.92"
$ws.Range("E36").Value = "  -1.89%  "

Set-TextValue "D37" "5.19"
$ws.Range("E37").Value = "  -5.72%  "

Set-TextValue "D38" "0.0411"
$ws.Range("E38").Value = "  -0.99%  "

Set-TextValue "D39" "0.0793"
$ws.Range("E39").Value = "  -2.37%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.048.64"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.120"
$ws.Range("E41").Value = "  -0.29%  "

Set-TextValue "D42" "2.69"
$ws.Range("E42").Value = "  -9.02%  "

$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("E44").Value = "  +1.42%  "

Set-TextValue "D46" "2.05"
$ws.Range("E46").Value = "  -3.96%  "

Set-TextValue "D47" "121.35"
$ws.Range("E47").Value = "  +1.62%  "

Set-TextValue "D48" "23.99"
$ws.Range("E48").Value = "  -5.78%  "

$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0504"
$ws.Range("E50").Value = "  -6.37%  "

$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D51" "2.38"
$ws.Range("E51").Value = "  +55.68%  "
